# FIMA_APP.xlsx — "added country brief info"
#
# The raw data-generation script was re-run, which reshuffled the random
# draws used to populate the "Interventions" (sheet "Interventions") and
# "Instruments" worksheets for the 2nd and 3rd countries (Aurelia / id 33,
# and Xenon / id 35) so that, like the first country (Ruritania / id 29),
# the intervention/instrument names now line up in ascending catalogue
# order within each type block. It also leaves the UI focused on the
# "Interventions" sheet (where the new data lives) instead of "About".

$wb  = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Interventions")
$ws3 = $wb.Worksheets.Item("Instruments")

# ---------------------------------------------------------------------
# 1. Interventions sheet (B17:B46) — re-sort each 8-row "Land Use" block
#    and 7-row "Protection Gap" block into ascending catalogue order.
# ---------------------------------------------------------------------
$landUse = @(
    "Silvopasture",
    "Reduced-till farming",
    "Dams and seawalls",
    "Restoring degraded forest",
    "Precision agriculture",
    "Agroforestry",
    "Large and medium scale irrigation",
    "Climate-resilient seeds"
)

$protectionGap = @(
    "Catastrophe bonds",
    "Insurance premium subsidies",
    "Microinsurance",
    "Cross-border reinsurance",
    "Compulsory insurance coverage",
    "Insurance bundling",
    "Risk-based solvency capital requirements"
)

# Aurelia block: rows 17-24 (Land Use), 25-31 (Protection Gap)
$row = 17
foreach ($val in $landUse) {
    $ws2.Cells.Item($row, 2).Value = $val
    $row = $row + 1
}
foreach ($val in $protectionGap) {
    $ws2.Cells.Item($row, 2).Value = $val
    $row = $row + 1
}

# Xenon block: rows 32-39 (Land Use), 40-46 (Protection Gap)
$row = 32
foreach ($val in $landUse) {
    $ws2.Cells.Item($row, 2).Value = $val
    $row = $row + 1
}
foreach ($val in $protectionGap) {
    $ws2.Cells.Item($row, 2).Value = $val
    $row = $row + 1
}

# ---------------------------------------------------------------------
# 2. Instruments sheet (B8:B19) — same re-sort for the Aurelia / Xenon
#    blocks (6 rows each).
# ---------------------------------------------------------------------
$instruments = @(
    "Sustainability-linked bonds",
    "Sustainability-linked loans",
    "Debt-for-nature swaps",
    "Carbon credits",
    "Biodiversity credits",
    "Credit enhancement"
)

# Aurelia block: rows 8-13
$row = 8
foreach ($val in $instruments) {
    $ws3.Cells.Item($row, 2).Value = $val
    $row = $row + 1
}

# Xenon block: rows 14-19
$row = 14
foreach ($val in $instruments) {
    $ws3.Cells.Item($row, 2).Value = $val
    $row = $row + 1
}

# ---------------------------------------------------------------------
# 3. View state — leave the Instruments sheet's selection on D18, then
#    finish with the Interventions sheet active/selected on B32:C46
#    (matches the saved workbook view in the edit).
# ---------------------------------------------------------------------
$ws3.Range("D18").Select()
$ws2.Activate()
$ws2.Range("B32:C46").Select()

try {
    $excel.ActiveWindow.ScrollRow = 16
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
    # Scroll position isn't part of the saved cell/selection state in
    # every host; ignore if unsupported.
}
